# Updates eligible population estimates for a number of projects
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with revised population figures / status notes ---

# Row 22 - Wolverhampton-ish entry: pop revised, status note replaced, fill set to white (theme)
$ws.Range("E22").Value = 44257
$ws.Range("E22").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1
$ws.Range("F22").Value = "06/03/2023 pop from Tim Windle (previously 30502). Pop needs confirming"

# Row 23 - pop figure confirmed (value unchanged) but fill/status updated
$ws.Range("E23").Value = 3682
$ws.Range("E23").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1
$ws.Range("F23").Value = "06/03/2023 - pop confirmed by Tim Windle. Pop needs confirming"

# Row 26 - pop revised, status note replaced, fill set to white
$ws.Range("E26").Value = 33033
$ws.Range("E26").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1
$ws.Range("F26").Value = "06/03/2023 - pop from Tim Windle (previously 130053). Pop needs confirming"

# Row 36 - pop confirmed (value unchanged) but fill/status updated
$ws.Range("E36").Value = 25662
$ws.Range("E36").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1
$ws.Range("F36").Value = "06/03/2023 pop confirmed by Tim Windle. figure from project - to check with Poppy"

# Row 37 - pop revised, status note replaced, fill set to white
$ws.Range("E37").Value = 51258
$ws.Range("E37").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1
$ws.Range("F37").Value = "06/03/2023 pop from Tim Windle (previously 142384) population advised from project (2023-02-17), to check with Poppy"

# Row 38 - eligible_population was blank, now populated; status note replaced
$ws.Range("E38").Value = 21315
$ws.Range("F38").Value = "06/03/2023 pop supplied by Tim Windle"

# --- Append four new project rows to the table (rows 39-42) ---

$tbl = $ws.ListObjects.Item(1)

$tbl.ListRows.Add() | Out-Null
$ws.Range("A39").Value = "52R"
$ws.Range("B39").Value = "Nottingham and Nottinghamshire"
$ws.Range("C39").Value = "Phase 3"
$ws.Range("D39").Value = 44927
$ws.Range("D39").NumberFormat = "mmm-yy"
$ws.Range("E39").Value = 32118
$ws.Range("F39").Value = "06/03/2023 pop supplied by Tim Windle"
$ws.Range("G39").Value = "E56000024"
$ws.Range("H39").Value = "East Midlands"

$tbl.ListRows.Add() | Out-Null
$ws.Range("A40").Value = "11J"
$ws.Range("B40").Value = "Dorset"
$ws.Range("C40").Value = "Phase 3"
$ws.Range("D40").Value = 44927
$ws.Range("D40").NumberFormat = "mmm-yy"
$ws.Range("E40").Value = 5660
$ws.Range("F40").Value = "06/03/2023 pop supplied by Tim Windle"
$ws.Range("G40").Value = "E56000016"
$ws.Range("H40").Value = "Wessex"

$tbl.ListRows.Add() | Out-Null
$ws.Range("A41").Value = "03H03K"
$ws.Range("B41").Value = "North and NE Lincolnshire"
$ws.Range("C41").Value = "Phase 3"
$ws.Range("D41").Value = 44927
$ws.Range("D41").NumberFormat = "mmm-yy"
$ws.Range("E41").Value = 45909
$ws.Range("F41").Value = "06/03/2023 pop supplied by Tim Windle"
$ws.Range("G41").Value = "E56000026"
$ws.Range("H41").Value = "Humber, Coast and Vale"

$tbl.ListRows.Add() | Out-Null
$ws.Range("A42").Value = "D2P2L"
$ws.Range("B42").Value = "Sandwell and West Birmingham"
$ws.Range("C42").Value = "Phase 3"
$ws.Range("D42").Value = 44927
$ws.Range("D42").NumberFormat = "mmm-yy"
$ws.Range("E42").Value = 34310
$ws.Range("F42").Value = "06/03/2023 pop supplied by Tim Windle"
$ws.Range("G42").Value = "E56000007"
$ws.Range("H42").Value = "West Midlands"

# Match the last active selection recorded in the saved workbook
$ws.Range("F42").Select() | Out-Null
